$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E2").Value = 21054
$ws.Range("F2").Value = 15651
$ws.Range("G2").Value = 3215
$ws.Range("H2").Value = 696
$ws.Range("J2").Value = 1343
$ws.Range("E5").Value = 90356
$ws.Range("F5").Value = 68024
$ws.Range("G5").Value = 10918
$ws.Range("H5").Value = 3957
$ws.Range("I5").Value = 1511
$ws.Range("J5").Value = 5946
$ws.Range("E6").Value = 996634
$ws.Range("F6").Value = 770481
$ws.Range("G6").Value = 152960
$ws.Range("H6").Value = 29922
$ws.Range("I6").Value = 7385
$ws.Range("J6").Value = 35886
$ws.Range("E7").Value = 190369
$ws.Range("F7").Value = 128924
$ws.Range("G7").Value = 18822
$ws.Range("H7").Value = 8638
$ws.Range("I7").Value = 6087
$ws.Range("J7").Value = 27898
$ws.Range("F8").Value = 33924
$ws.Range("G8").Value = 4669
$ws.Range("H8").Value = 2485
$ws.Range("I8").Value = 2809
$ws.Range("J8").Value = 15749
$ws.Range("E9").Value = 774599
$ws.Range("F9").Value = 426262
$ws.Range("G9").Value = 71245
$ws.Range("H9").Value = 42263
$ws.Range("I9").Value = 31177
$ws.Range("J9").Value = 203652
$ws.Range("E10").Value = 12607
$ws.Range("F10").Value = 5572
$ws.Range("G10").Value = 895
$ws.Range("H10").Value = 865
$ws.Range("I10").Value = 855
$ws.Range("J10").Value = 4420
$ws.Range("E12").Value = 636398
$ws.Range("F12").Value = 482049
$ws.Range("G12").Value = 106661
$ws.Range("H12").Value = 29918
$ws.Range("I12").Value = 4846
$ws.Range("J12").Value = 12924
$ws.Range("E13").Value = 143481
$ws.Range("G13").Value = 19816
$ws.Range("E14").Value = 76782
$ws.Range("F14").Value = 54769
$ws.Range("G14").Value = 9150
$ws.Range("H14").Value = 3704
$ws.Range("J14").Value = 7408
$ws.Range("E15").Value = 55369
$ws.Range("F15").Value = 34520
$ws.Range("G15").Value = 7774
$ws.Range("H15").Value = 3069
$ws.Range("I15").Value = 1684
$ws.Range("J15").Value = 8322
$ws.Range("E16").Value = 29631
$ws.Range("F16").Value = 19413
$ws.Range("G16").Value = 4311
$ws.Range("H16").Value = 1736
$ws.Range("J16").Value = 3223
$ws.Range("E18").Value = 189351
$ws.Range("F18").Value = 145953
$ws.Range("G18").Value = 28312
$ws.Range("H18").Value = 7691
$ws.Range("I18").Value = 1752
$ws.Range("J18").Value = 5643
$ws.Range("E19").Value = 277543
$ws.Range("F19").Value = 106098
$ws.Range("G19").Value = 20840
$ws.Range("H19").Value = 16268
$ws.Range("I19").Value = 18994
$ws.Range("J19").Value = 115343
$ws.Range("E21").Value = 12520
$ws.Range("F21").Value = 5980
$ws.Range("G21").Value = 1517
$ws.Range("I21").Value = 934
$ws.Range("E22").Value = 379657
$ws.Range("F22").Value = 247760
$ws.Range("G22").Value = 36036
$ws.Range("H22").Value = 20495
$ws.Range("I22").Value = 12266
$ws.Range("J22").Value = 63100
$ws.Range("E23").Value = 36362
$ws.Range("F23").Value = 24051
$ws.Range("G23").Value = 4749
$ws.Range("H23").Value = 1995
$ws.Range("I23").Value = 967
$ws.Range("J23").Value = 4600
$ws.Range("E24").Value = 153717
$ws.Range("F24").Value = 80620
$ws.Range("G24").Value = 15417
$ws.Range("H24").Value = 9184
$ws.Range("I24").Value = 7748
$ws.Range("J24").Value = 40748
$ws.Range("E26").Value = 2954
$ws.Range("G26").Value = 449
$ws.Range("H26").Value = 229
$ws.Range("I26").Value = 69
$ws.Range("J26").Value = 759
$ws.Range("F29").Value = 521
$ws.Range("G29").Value = 99
$ws.Range("H29").Value = 110
$ws.Range("I29").Value = 77
$ws.Range("J29").Value = 731
$ws.Range("E32").Value = 1224
$ws.Range("F32").Value = 382
$ws.Range("G32").Value = 69
$ws.Range("H32").Value = 49
$ws.Range("I32").Value = 29
$ws.Range("E33").Value = 20529
$ws.Range("F33").Value = 12683
$ws.Range("G33").Value = 2945
$ws.Range("H33").Value = 1819
$ws.Range("J33").Value = 2466
$ws.Range("E35").Value = 2463
$ws.Range("F35").Value = 1302
$ws.Range("G35").Value = 280
$ws.Range("H35").Value = 200
$ws.Range("I35").Value = 110
$ws.Range("J35").Value = 571
$ws.Range("E37").Value = 879
$ws.Range("F37").Value = 645
$ws.Range("G37").Value = 96
